# "Initial Schematic, Item list update"
# Populates the ShowHerb item list (rows 2-6) on the Accounting sheet with
# the parts used for the project, and leaves the selection on G6 as the
# last-used cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Arduino Nano 33 IoT
$ws.Range("B2").Value = "Arduino Nano 33 IoT "
$ws.Range("C2").Value = "Amazon"
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 15
$ws.Range("G2").Value = "https://www.amazon.com/gp/product/B07VW9TSKD/ref=ppx_yo_dt_b_asin_title_o01_s00?ie=UTF8&psc=1"

# Row 3 - Adafruit I2C Multiplexer
$ws.Range("B3").Value = "Adafruit I2C Multiplexer"
$ws.Range("C3").Value = "Amazon"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 8.5
$ws.Range("G3").Value = "https://www.amazon.com/gp/product/B017C09ETS/ref=ppx_yo_dt_b_asin_title_o00_s00?ie=UTF8&psc=1"

# Row 4 - STEMMA Soil Sensor
$ws.Range("B4").Value = "STEMMA Soil Sensor "
$ws.Range("C4").Value = "Adafruit"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 7.5
$ws.Range("G4").Value = "https://www.adafruit.com/product/4026?gclid=Cj0KCQjwl_SHBhCQARIsAFIFRVXqnetQca0YiH1L4WKDv7rMAurkYQ8s318Mtg1VIUXXEpn5wohOilwaAkpSEALw_wcB"

# Row 5 - 12V 1/2'' Electric Solenoid Valve
$ws.Range("B5").Value = "12V 1/2'' Electric Solenoid Valve"
$ws.Range("C5").Value = "Amazon"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 13
$ws.Range("G5").Value = "https://www.amazon.com/gp/product/B07NWCRM75/ref=ppx_yo_dt_b_asin_title_o01_s00?ie=UTF8&psc=1"

# Row 6 - 5V Relay / Elegoo Kit
$ws.Range("B6").Value = "5V Relay "
$ws.Range("C6").Value = "Elegoo Kit "
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0

# Match the saved selection state from the edit (active cell G6)
$ws.Range("G6").Select()
